$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) C28: trim the trailing "inserisci nuovo poi" quoted clause
$ws.Range("C28").Value = 'si inizia lo sviluppo dei casi d''uso "inserisci nuovo comune" '

# 2) C35: "altri" -> "tutti"
$ws.Range("C35").Value = "Si inseriscono paragrafi descrittivi brevi su tutti i casi d'uso "

# 3) Stash the two source styles we'll need for the rebuilt block 38 onto
#    scratch cells, since their original rows are about to be deleted.
$ws.Range("A24").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("C45:H45").Copy()
$ws.Range("Z2:AE2").PasteSpecial(-4122)

# 4) Remove the old rows 36-50 block (will be rebuilt below), keeping row 52's
#    content which shifts up to become the new final row.
$ws.Rows("36:50").Delete()

# 5) Re-open a gap of 8 rows so the old row 52 (now row 37) becomes row 45 again,
#    leaving rows 36-44 blank for the new content.
$ws.Rows("37:44").Insert()

# 6) Row 38: date-block header row (A38 "dettaglio poi", B38 = iterazione 2,
#    C38:H38 merged with the "Si identificano..." paragraph, style copied from
#    the old merged banner row).
$ws.Range("Z1").Copy()
$ws.Range("A38").PasteSpecial(-4122)
$ws.Range("A38").Value = "dettaglio poi"

$ws.Range("B34").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("B38").Value = 2

$ws.Range("Z2:AE2").Copy()
$ws.Range("C38:H38").PasteSpecial(-4122)
$ws.Range("C38:H38").Merge()
$ws.Range("C38").Value = "Si identificano altri 5 casi d'uso pari a circa il 15 % del totale e si provvede alla definizione della versione dettagliata"

# clear the scratch cells
$ws.Range("Z1:AE2").Clear()

# 6) Rows 39-44: shaded "C" cells (style copied from the existing shaded block)
#    plus plain "ok" notes in column D.
$ws.Range("C25").Copy()
$ws.Range("C39:C44").PasteSpecial(-4122)

$ws.Range("C39").Value = "valida contenuto"
$ws.Range("D39").Value = "ok"

$ws.Range("C40").Value = "Registrazione nuovo utente"
$ws.Range("D40").Value = "ok "

$ws.Range("C41").Value = "autocertificazione contento"
$ws.Range("D41").Value = "ok"

$ws.Range("C42").Value = "registrazione nuovo utente"
$ws.Range("D42").Value = "ok"

$ws.Range("C43").Value = "autenticazione utente"
$ws.Range("D43").Value = "ok"

$ws.Range("C44").Value = "assegna ruolo utente"
$ws.Range("D44").Value = "ok"

# 7) View/window cosmetics to match the saved workbook state.
$ws.Range("E10").Select()
$excel.ActiveWindow.Zoom = 149
$excel.ActiveWindow.ScrollRow = 34
